$d = $word.ActiveDocument

function Merge-RequirementId {
    param(
        [int]$ParaIndex,
        [string]$IdText
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $r = $p.Range
    $null = $r.Find.Execute($IdText, $false, $false, $false, $false, $false, $true, 1, $false, $IdText, 2)
}

# SSS-0002 (paragraph 4): merge the split "O sistema DE" + "VE " runs into one run.
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$null = $r4.Find.Execute("O sistema DEVE ", $false, $false, $false, $false, $false, $true, 1, $false, "O sistema DEVE ", 2)

# SSS-0005 (paragraph 7): merge the split "SSS-000" + "5: " runs into one run ...
Merge-RequirementId 7 "SSS-0005: "
# ... and drop "e validação " from the body text.
$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
$null = $r7.Find.Execute("a visualização e validação dos", $false, $false, $false, $false, $false, $true, 1, $false, "a visualização dos", 2)

# SSS-0007 .. SSS-0017 (paragraphs 9 .. 19): merge the split bold requirement-id
# runs (e.g. "SSS-000" + "7: ") into a single run, leaving the body text as-is.
Merge-RequirementId 9  "SSS-0007: "
Merge-RequirementId 10 "SSS-0008: "
Merge-RequirementId 11 "SSS-0009: "
Merge-RequirementId 12 "SSS-0010: "
Merge-RequirementId 13 "SSS-0011: "
Merge-RequirementId 14 "SSS-0012: "
Merge-RequirementId 15 "SSS-0013: "
Merge-RequirementId 16 "SSS-0014: "
Merge-RequirementId 17 "SSS-0015: "
Merge-RequirementId 18 "SSS-0016: "
Merge-RequirementId 19 "SSS-0017: "
